# Update header and append pence values to "Price Per Unit" column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Price Per Unit (£s)"

$values = @(
    1054.63,
    1629.57,
    1606.67,
    1201.03,
    1434.28,
    289.95,
    500.08,
    2100.87,
    132.42,
    994.4299999999999,
    896.17,
    905.11,
    1076.55,
    1592.19,
    325.24,
    1058.77,
    291.21,
    1255.21,
    2135.78,
    497.25,
    1129.05,
    854.6900000000001,
    830.78,
    1131.27,
    2449.87,
    1578.02,
    2297.42,
    2606.77,
    277.81,
    2172.67,
    2072.84,
    888.53,
    618.09,
    644.99,
    2075.84,
    2581.93,
    724.99,
    276.75,
    1432.77,
    734.76,
    39.34,
    2900.69,
    1199.25,
    847.83,
    650.3099999999999,
    1495.55,
    2106.64,
    2213.31,
    1038.88,
    1595.77
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
